$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299.75
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 899
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = -1125
$ws.Range("H33").Value = 532.4
$ws.Range("J33").Value = 850
$ws.Range("L33").Value = 850
$ws.Range("N33").Value = -1308
$ws.Range("H43").Value = 3935.05
$ws.Range("J43").Value = 2970.5557
$ws.Range("L43").Value = 2970.5557
$ws.Range("N43").Value = -3108.5557
$ws.Range("H69").Value = 13909.637
$ws.Range("J69").Value = 17500
$ws.Range("L69").Value = 52500
$ws.Range("N69").Value = -54248
$ws.Range("H72").Value = 13909.637
$ws.Range("J72").Value = 17500
$ws.Range("L72").Value = 157500
$ws.Range("N72").Value = -166236
$ws.Range("H116").Value = 13207.733
$ws.Range("I116").Value = 5167.778
$ws.Range("J116").Value = 25267.666
$ws.Range("K116").Value = 5167.778
$ws.Range("L116").Value = 25267.666
$ws.Range("M116").Value = -1725.778
$ws.Range("N116").Value = -32151.666
$ws.Range("H125").Value = 15626
$ws.Range("I125").Value = 20430
$ws.Range("J125").Value = 13224
$ws.Range("K125").Value = 183870
$ws.Range("L125").Value = 119016
$ws.Range("M125").Value = -181410
$ws.Range("N125").Value = -123936
$ws.Range("H131").Value = 2274646
$ws.Range("I131").Value = 2132.111
$ws.Range("K131").Value = 6396.333
$ws.Range("M131").Value = -1356.333
$ws.Range("H138").Value = 3413.3394
$ws.Range("I138").Value = 1866.8064
$ws.Range("K138").Value = 5600.4192
$ws.Range("M138").Value = -460.4192000000003

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 28332.5
$ws.Range("J43").Value = 28332.5
$ws.Range("L43").Value = 28332.5
$ws.Range("N43").Value = -28958.5
$ws.Range("H63").Value = 4665.3335
$ws.Range("I63").Value = 4665.3335
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4665.3335
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3979.3335
$ws.Range("H66").Value = 4665.3335
$ws.Range("I66").Value = 4665.3335
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 23326.6675
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -19894.6675
$ws.Range("H110").Value = 5712.3
$ws.Range("I110").Value = 5416.6665
$ws.Range("J110").Value = 6599.2
$ws.Range("K110").Value = 5416.6665
$ws.Range("L110").Value = 6599.2
$ws.Range("M110").Value = -3371.6665
$ws.Range("N110").Value = -10689.2
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 99998.336
$ws.Range("J125").Value = 99998.336
$ws.Range("L125").Value = 99998.336
$ws.Range("N125").Value = -109838.336
$ws.Range("H129").Value = 60699
$ws.Range("J129").Value = 60699
$ws.Range("L129").Value = 60699
$ws.Range("N129").Value = -70699
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079
$ws.Range("N130").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27782424
$ws.Range("I31").Value = 76927050
$ws.Range("J31").Value = 5030.2607
$ws.Range("K31").Value = 76927050
$ws.Range("L31").Value = 5030.2607
$ws.Range("M31").Value = -76926755
$ws.Range("N31").Value = -5620.2607
$ws.Range("H34").Value = 27782424
$ws.Range("I34").Value = 76927050
$ws.Range("J34").Value = 5030.2607
$ws.Range("K34").Value = 76927050
$ws.Range("L34").Value = 5030.2607
$ws.Range("M34").Value = -76926848
$ws.Range("N34").Value = -5434.2607
$ws.Range("H58").Value = 2115.6843
$ws.Range("I58").Value = 1932.4073
$ws.Range("J58").Value = 2565.5454
$ws.Range("K58").Value = 1932.4073
$ws.Range("L58").Value = 2565.5454
$ws.Range("M58").Value = -1729.4073
$ws.Range("N58").Value = -2971.5454
$ws.Range("H136").Value = 2115.6843
$ws.Range("I136").Value = 1932.4073
$ws.Range("J136").Value = 2565.5454
$ws.Range("K136").Value = 5797.2219
$ws.Range("L136").Value = 7696.6362
$ws.Range("M136").Value = -3247.2219
$ws.Range("N136").Value = -12796.6362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 6247.4287
$ws.Range("I95").Value = 1733.1666
$ws.Range("J95").Value = 33333
$ws.Range("K95").Value = 5199.4998
$ws.Range("L95").Value = 99999
$ws.Range("M95").Value = -3140.4998
$ws.Range("N95").Value = -104117
$ws.Range("H116").Value = 6498.2856
$ws.Range("I116").Value = 2025.8334
$ws.Range("K116").Value = 6077.5002
$ws.Range("M116").Value = -2635.5002
$ws.Range("H122").Value = 55466.668
$ws.Range("I122").Value = 66360
$ws.Range("K122").Value = 597240
$ws.Range("M122").Value = -594790

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4279.4287
$ws.Range("I97").Value = 692.2727
$ws.Range("J97").Value = 17432.334
$ws.Range("K97").Value = 692.2727
$ws.Range("L97").Value = 17432.334
$ws.Range("M97").Value = -196.2727
$ws.Range("N97").Value = -18424.334
$ws.Range("H102").Value = 1429.871
$ws.Range("I102").Value = 1356.069
$ws.Range("K102").Value = 1356.069
$ws.Range("M102").Value = 265.931
$ws.Range("H113").Value = 1687063.1
$ws.Range("I113").Value = 3784
$ws.Range("K113").Value = 3784
$ws.Range("M113").Value = -1614
$ws.Range("H122").Value = 2541844.8
$ws.Range("I122").Value = 3302320.5
$ws.Range("J122").Value = 6925.3335
$ws.Range("K122").Value = 9906961.5
$ws.Range("L122").Value = 20776.0005
$ws.Range("M122").Value = -9904511.5
$ws.Range("N122").Value = -25676.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7825717
$ws.Range("I100").Value = 6000
$ws.Range("J100").Value = 8077966
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 8077966
$ws.Range("M100").Value = -5459
$ws.Range("N100").Value = -8079048
$ws.Range("H132").Value = 3356
$ws.Range("I132").Value = 1993.75
$ws.Range("K132").Value = 5981.25
$ws.Range("M132").Value = -3451.25
$ws.Range("H136").Value = 8421.777
$ws.Range("J136").Value = 10078
$ws.Range("L136").Value = 30234
$ws.Range("N136").Value = -35334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H122").Value = 1934.62
$ws.Range("I122").Value = 1743.1578
$ws.Range("J122").Value = 2540.9167
$ws.Range("K122").Value = 5229.4734
$ws.Range("L122").Value = 7622.750100000001
$ws.Range("M122").Value = -2779.4734
$ws.Range("N122").Value = -12522.7501
$ws.Range("H126").Value = 2753.8057
$ws.Range("I126").Value = 1940.6428
$ws.Range("J126").Value = 5599.875
$ws.Range("K126").Value = 5821.928400000001
$ws.Range("L126").Value = 16799.625
$ws.Range("M126").Value = -3351.928400000001
$ws.Range("N126").Value = -21739.625
$ws.Range("H132").Value = 2921.6924
$ws.Range("I132").Value = 3254.5715
$ws.Range("K132").Value = 9763.7145
$ws.Range("M132").Value = -7233.7145
$ws.Range("N50").ClearContents()
